$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44452

# Row 3
$ws.Range("D3").Value = 44461
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 31000
$ws.Range("O3").Value = 32000
$ws.Range("P3").Value = 31500
$ws.Range("S3").Value = 3150

# Row 4
$ws.Range("D4").Value = 44461
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 30000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 30000
$ws.Range("S4").Value = 3000

# Row 5
$ws.Range("D5").Value = 44447
$ws.Range("L5").Value = "Primera"
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 2150

# Row 6
$ws.Range("D6").Value = 44487
$ws.Range("N6").Value = 23000
$ws.Range("O6").Value = 24000
$ws.Range("P6").Value = 23500
$ws.Range("S6").Value = 2350

# Row 10
$ws.Range("D10").Value = 44446
